# Auto-applied data refresh for Sheets/Maduin_Profits.xlsx (per scheduled runner diff).
# Updates currentAveragePrice/NQ/HQ, LevePriceNQ/HQ and LeveProfitNQ/HQ columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1329.4166
$ws.Range("I6").Value = 635.3333
$ws.Range("K6").Value = 1905.9999
$ws.Range("M6").Value = -1793.9999
$ws.Range("H40").Value = 2299.3333
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2299.3333
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2649.3333
$ws.Range("H52").Value = 3000
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H58").Value = 6671.6665
$ws.Range("I58").Value = 15
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 45
$ws.Range("L58").Value = 30000
$ws.Range("M58").Value = 105
$ws.Range("N58").Value = -30300
$ws.Range("H88").Value = 3458
$ws.Range("I88").Value = 3998
$ws.Range("J88").Value = 3350
$ws.Range("K88").Value = 3998
$ws.Range("L88").Value = 3350
$ws.Range("M88").Value = -3592
$ws.Range("N88").Value = -4162
$ws.Range("H91").Value = 3458
$ws.Range("I91").Value = 3998
$ws.Range("J91").Value = 3350
$ws.Range("K91").Value = 3998
$ws.Range("L91").Value = 3350
$ws.Range("M91").Value = -2594
$ws.Range("N91").Value = -6158
$ws.Range("H138").Value = 2919.7778
$ws.Range("I138").Value = 2611.1428
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 7833.428400000001
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = -2693.428400000001
$ws.Range("N138").Value = -22280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3309.3333
$ws.Range("I88").Value = 1336.5
$ws.Range("J88").Value = 4098.467
$ws.Range("K88").Value = 1336.5
$ws.Range("L88").Value = 4098.467
$ws.Range("M88").Value = -930.5
$ws.Range("N88").Value = -4910.467
$ws.Range("H91").Value = 3309.3333
$ws.Range("I91").Value = 1336.5
$ws.Range("J91").Value = 4098.467
$ws.Range("K91").Value = 1336.5
$ws.Range("L91").Value = 4098.467
$ws.Range("M91").Value = 67.5
$ws.Range("N91").Value = -6906.467
$ws.Range("H110").Value = 970
$ws.Range("I110").Value = 970
$ws.Range("K110").Value = 970
$ws.Range("M110").Value = 1075
$ws.Range("H132").Value = 1121.375
$ws.Range("I132").Value = 1121.375
$ws.Range("K132").Value = 3364.125
$ws.Range("M132").Value = -834.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3476.3125
$ws.Range("I105").Value = 3086.3076
$ws.Range("K105").Value = 3086.3076
$ws.Range("M105").Value = -1339.3076

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 969
$ws.Range("I16").Value = 700.5
$ws.Range("K16").Value = 700.5
$ws.Range("M16").Value = -413.5
$ws.Range("H31").Value = 2228
$ws.Range("J31").Value = 2499
$ws.Range("L31").Value = 2499
$ws.Range("N31").Value = -3089
$ws.Range("H34").Value = 2228
$ws.Range("J34").Value = 2499
$ws.Range("L34").Value = 2499
$ws.Range("N34").Value = -2903
$ws.Range("H58").Value = 3556
$ws.Range("I58").Value = 2112
$ws.Range("K58").Value = 2112
$ws.Range("M58").Value = -1909
$ws.Range("H105").Value = 1388.5714
$ws.Range("I105").Value = 1433.6666
$ws.Range("K105").Value = 1433.6666
$ws.Range("M105").Value = 313.3334
$ws.Range("H107").Value = 966.6667
$ws.Range("I107").Value = 950
$ws.Range("K107").Value = 950
$ws.Range("M107").Value = 970
$ws.Range("H113").Value = 969
$ws.Range("I113").Value = 700.5
$ws.Range("K113").Value = 700.5
$ws.Range("M113").Value = 1469.5
$ws.Range("H134").Value = 2346.1333
$ws.Range("I134").Value = 1827.1111
$ws.Range("K134").Value = 5481.3333
$ws.Range("M134").Value = -2946.3333
$ws.Range("H136").Value = 3556
$ws.Range("I136").Value = 2112
$ws.Range("K136").Value = 6336
$ws.Range("M136").Value = -3786

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 267.6
$ws.Range("I7").Value = 86.666664
$ws.Range("K7").Value = 259.999992
$ws.Range("M7").Value = -147.999992
$ws.Range("H17").Value = 834.1
$ws.Range("I17").Value = 450
$ws.Range("J17").Value = 876.7778
$ws.Range("K17").Value = 1350
$ws.Range("L17").Value = 2630.3334
$ws.Range("M17").Value = -1181
$ws.Range("N17").Value = -2968.3334
$ws.Range("H26").Value = 800
$ws.Range("J26").Value = 800
$ws.Range("L26").Value = 2400
$ws.Range("N26").Value = -2976
$ws.Range("H55").Value = 5251
$ws.Range("H80").Value = 6847.8
$ws.Range("I80").Value = 3079.8333
$ws.Range("J80").Value = 12499.75
$ws.Range("K80").Value = 9239.499899999999
$ws.Range("L80").Value = 37499.25
$ws.Range("M80").Value = -8303.499899999999
$ws.Range("N80").Value = -39371.25
$ws.Range("H83").Value = 6847.8
$ws.Range("I83").Value = 3079.8333
$ws.Range("J83").Value = 12499.75
$ws.Range("K83").Value = 27718.4997
$ws.Range("L83").Value = 112497.75
$ws.Range("M83").Value = -23038.4997
$ws.Range("N83").Value = -121857.75
$ws.Range("H92").Value = 574.75
$ws.Range("I92").Value = 274.5
$ws.Range("K92").Value = 823.5
$ws.Range("M92").Value = 424.5
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H131").Value = 1037.1666
$ws.Range("J131").Value = 1094.8462
$ws.Range("L131").Value = 3284.5386
$ws.Range("N131").Value = -13364.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1677.1
$ws.Range("I113").Value = 1603.25
$ws.Range("K113").Value = 1603.25
$ws.Range("M113").Value = 566.75
$ws.Range("H126").Value = 4457.8
$ws.Range("I126").Value = 4457.8
$ws.Range("K126").Value = 13373.4
$ws.Range("M126").Value = -10903.4
$ws.Range("H132").Value = 7168
$ws.Range("I132").Value = 8322.666999999999
$ws.Range("K132").Value = 24968.001
$ws.Range("M132").Value = -22438.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2218.625
$ws.Range("I22").Value = 1203.8334
$ws.Range("J22").Value = 5263
$ws.Range("K22").Value = 1203.8334
$ws.Range("L22").Value = 5263
$ws.Range("M22").Value = -908.8334
$ws.Range("N22").Value = -5853
$ws.Range("H27").Value = 2218.625
$ws.Range("I27").Value = 1203.8334
$ws.Range("J27").Value = 5263
$ws.Range("K27").Value = 1203.8334
$ws.Range("L27").Value = 5263
$ws.Range("M27").Value = -1096.8334
$ws.Range("N27").Value = -5477
$ws.Range("H46").Value = 1479.5
$ws.Range("I46").Value = 1263
$ws.Range("K46").Value = 1263
$ws.Range("M46").Value = -1075
$ws.Range("H61").Value = 2000
$ws.Range("J61").Value = 2000
$ws.Range("L61").Value = 2000
$ws.Range("N61").Value = -2404
$ws.Range("H100").Value = 300
$ws.Range("I100").Value = 300
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 300
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 2000
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -6340
$ws.Range("H132").Value = 1932.75
$ws.Range("I132").Value = 1998.8
$ws.Range("J132").Value = 1602.5
$ws.Range("K132").Value = 5996.4
$ws.Range("L132").Value = 4807.5
$ws.Range("M132").Value = -3466.4
$ws.Range("N132").Value = -9867.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 908.8
$ws.Range("I96").Value = 965.3333
$ws.Range("J96").Value = 400
$ws.Range("K96").Value = 965.3333
$ws.Range("L96").Value = 400
$ws.Range("M96").Value = 407.6667
$ws.Range("N96").Value = -3146
$ws.Range("H113").Value = 760
$ws.Range("I113").Value = 720
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 2160
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = 10
$ws.Range("N113").Value = -6740
$ws.Range("H126").Value = 3213.4
$ws.Range("I126").Value = 3226
$ws.Range("K126").Value = 9678
$ws.Range("M126").Value = -7208
$ws.Range("H132").Value = 2552.3572
$ws.Range("I132").Value = 2517.5386
$ws.Range("J132").Value = 3005
$ws.Range("K132").Value = 7552.6158
$ws.Range("L132").Value = 9015
$ws.Range("M132").Value = -5022.6158
$ws.Range("N132").Value = -14075
